$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column values for rows 8-16
$ws.Range("H8").Value = 15.1610087863751
$ws.Range("H9").Value = 11.81293910030774
$ws.Range("H10").Value = 10.85437886210314
$ws.Range("H11").Value = 14.83900377952033

# Row 12 updates
$ws.Range("B12").Value = 2500
$ws.Range("E12").Value = "Config_198"
$ws.Range("F12").Value = 45.45433982
$ws.Range("H12").Value = 12.74190179965575

$ws.Range("H13").Value = 11.19937228142982

# Row 14 updates
$ws.Range("D14").Value = 20
$ws.Range("E14").Value = "Config_152"
$ws.Range("F14").Value = 33.91581593
$ws.Range("H14").Value = 15.72664406951813

$ws.Range("H15").Value = 12.64418896608578

$ws.Range("H16").Value = 10.21370698986457

# Remove rows 17 and 18 entirely
$ws.Range("A17:H18").EntireRow.Delete()
